# Update the "modified time" stamp column on every portfolio sheet from
# 202509211530 -> 202509211531 (Web UI re-sync at 2025-09-21 07:31).

$wb = $excel.ActiveWorkbook

$oldStamp = "202509211530"
$newStamp = "202509211531"

# Sheet 1 "大智投资组合": timestamp lives in column E, data rows 2-9
$ws1 = $wb.Worksheets.Item(1)
for ($r = 2; $r -le 9; $r++) {
    $cell = $ws1.Cells.Item($r, 5)
    if ($cell.Value2 -eq $oldStamp) {
        $cell.Value = $newStamp
    }
}

# Sheet 2 "大成投资组合": timestamp lives in column E, data rows 2-11
$ws2 = $wb.Worksheets.Item(2)
for ($r = 2; $r -le 11; $r++) {
    $cell = $ws2.Cells.Item($r, 5)
    if ($cell.Value2 -eq $oldStamp) {
        $cell.Value = $newStamp
    }
}

# Sheet 3 "我的投资组合": timestamp lives in column G, data rows 2-13
$ws3 = $wb.Worksheets.Item(3)
for ($r = 2; $r -le 13; $r++) {
    $cell = $ws3.Cells.Item($r, 7)
    if ($cell.Value2 -eq $oldStamp) {
        $cell.Value = $newStamp
    }
}
